# Notas da Lista 02 e da Prova 01 adicionadas
#
# 1) Add the "Lista 03" header (D6) + the "Lista 02" column (C) of data to
#    the existing "Listas" sheet.
# 2) Create a new "Provas" sheet (as a copy of "Listas" so it inherits all
#    formatting/styles), positioned before "Listas", and rework its
#    content/layout into the "Provas" layout.

$wb = $excel.ActiveWorkbook
$wsListas = $wb.Worksheets.Item(1)
$wsListas.Name = "Listas"

# =====================================================================
# Step 1: add "Lista 03" header + "Lista 02" data to "Listas"
# =====================================================================

$wsListas.Range("D6").Value = "Lista 03"
$wsListas.Range("D6").HorizontalAlignment = -4108
$wsListas.Range("D6").Font.Bold = $true

# --- Step 2: duplicate "Listas" before itself, rename the copy "Provas" ---
$wsListas.Copy($wsListas)
$wsProvas = $wb.Worksheets.Item(1)
$wsProvas.Name = "Provas"
$wsListas = $wb.Worksheets.Item(2)

# =====================================================================
# Step 3: turn the "Provas" copy into the real Provas sheet
# =====================================================================

# Remove row 3 (Provas only has 2 note rows, not 3) and its merge.
$wsProvas.Range("A3:G3").UnMerge()
$wsProvas.Rows.Item(3).Clear()

# Provas header row (row 6): Matrícula / Prova 01 / Prova 02.  Clear the
# "Lista 03" header inherited from Listas' column D first.
$wsProvas.Range("D6").Clear()
$wsProvas.Range("A6").Value = "Matrícula"
$wsProvas.Range("B6").Value = "Prova 01"
$wsProvas.Range("A6:C6").HorizontalAlignment = -4108
$wsProvas.Range("A6:C6").Font.Bold = $true

# Extend the note rows from column G to column J.
$wsProvas.Range("A1:G1").UnMerge()
$wsProvas.Range("A2:G2").UnMerge()

$wsProvas.Range("A1").Value = "As notas das provas variam entre 0 e 100"
$wsProvas.Range("A2").Value = "A nota de cada prova será multiplicada por 0.25 para fins de computação da nota final da disciplina"

$wsProvas.Range("A1:J1").Merge()
$wsProvas.Range("A2:J2").Merge()

$wsProvas.Range("C6").Value = "Prova 02"
$wsProvas.Range("C6").HorizontalAlignment = -4108
$wsProvas.Range("C6").Font.Bold = $true

# Prova 01 grades (column B) as formulas; column C (Prova 02) is left blank,
# only the header was added.  Clear out the inherited Lista 01 / Lista 02
# grades first (B7:D26) before (re)writing column B with Prova 01 formulas.
$wsProvas.Range("B7:D26").ClearContents()

$wsProvas.Range("B7").Formula = "=(15+22+20+35)"
$wsProvas.Range("B8").Formula = "=(15+25+20+25)"
$wsProvas.Range("B9").Formula = "=(12+23+20+25)"
$wsProvas.Range("B10").Formula = "=(19+16+20+35)"
$wsProvas.Range("B11").Formula = "=(16+25+18+25)"
$wsProvas.Range("B12").Formula = "=(15+25+20+25)"
$wsProvas.Range("B13").Formula = "=(10+22+15+25)"
$wsProvas.Range("B14").Formula = "=(19+16+20+25)"
$wsProvas.Range("B15").Formula = "=(8+25+20+25)"
$wsProvas.Range("B16").Formula = "=(13+24+20+35)"
$wsProvas.Range("B17").Formula = "=(8+25+15+30)"
$wsProvas.Range("B18").Formula = "=(15+24+20+25)"
$wsProvas.Range("B19").Formula = "=(12+23+20+35)"
$wsProvas.Range("B20").Formula = "=(16+25+20+35)"
$wsProvas.Range("B21").Formula = "=(12+25+15+25)"
$wsProvas.Range("B22").Formula = "=(10+15+25+25)"
$wsProvas.Range("B23").Formula = "=(13+23+20+25)"
$wsProvas.Range("B24").Formula = "=(15+25+20+25)"
$wsProvas.Range("B25").Formula = "=(10+25+15+25)"
$wsProvas.Range("B26").Formula = "=(16+25+20+35)"

# =====================================================================
# Step 4: fill in "Lista 02" (column C) data on "Listas"
# =====================================================================

$wsListas.Range("C7").Value = 85
$wsListas.Range("C8").Value = 100
$wsListas.Range("C9").Value = 40
$wsListas.Range("C10").Value = 80
$wsListas.Range("C11").Value = 100
$wsListas.Range("C12").Value = 100
$wsListas.Range("C13").Value = 95
$wsListas.Range("C14").Value = 80
$wsListas.Range("C15").Value = 85
$wsListas.Range("C16").Value = 100
$wsListas.Range("C17").Value = 95
$wsListas.Range("C18").Value = 0
$wsListas.Range("C19").Value = 75
$wsListas.Range("C20").Value = 100
$wsListas.Range("C21").Value = 70
$wsListas.Range("C22").Value = 65
$wsListas.Range("C23").Value = 85
$wsListas.Range("C24").Value = 75
$wsListas.Range("C25").Value = 100
$wsListas.Range("C26").Value = 95

# Selection: "Listas" ends up with B31 selected (and not the active tab);
# "Provas" ends up active with A6 selected.  Select Listas first so that
# the final/active selection (and active tab) is on Provas.
$wsListas.Range("B31").Select()
$wsProvas.Activate()
$wsProvas.Range("A6").Select()

$wb.Save()
